$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6000
$ws.Range("B3").Value = 163075
$ws.Range("B4").Value = 30649
$ws.Range("B5").Value = 18.79
$ws.Range("B6").Value = 27.18
$ws.Range("B7").Value = 7.84
